$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source column formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values (coin rankings refreshed by the scraper)
$ws.Range("D2").Value = "63.883.63"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "3.351.47"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "547.96"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "172.74"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -3.60%  "
$ws.Range("D8").Value = "3.331.01"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "0.615"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "53.84"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "8.94"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "3.858.62"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "17.96"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.353.34"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").Value = "11.75"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "63.700.99"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "0.978"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "413.14"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "4.04"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "13.87"
$ws.Range("E24").Value = "  +14.97%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "4.33"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").Value = "83.14"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").Value = "10.59"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "2.73"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "8.63"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "29.16"
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").Value = "6.42"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "11.36"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("D33").Value = "578.24"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").Value = "57.99"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "0.148"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "35.30"
$ws.Range("E38").Value = "  -4.62%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.41"
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "0.0₃0742"
$ws.Range("E40").Value = "  -4.07%  "
$ws.Range("D41").Value = "0.368"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "3.151.60"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").Value = "0.0401"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -4.25%  "
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -4.64%  "
$ws.Range("D49").Value = "0.128"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").Value = "132.12"
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").Value = "  -3.41%  "
